$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1289.85
$ws.Range("I43").Value = 1565.5
$ws.Range("J43").Value = 1014.2
$ws.Range("K43").Value = 1565.5
$ws.Range("L43").Value = 1014.2
$ws.Range("M43").Value = -1496.5
$ws.Range("N43").Value = -1152.2

# Row 100
$ws.Range("H100").Value = 4707.5
$ws.Range("I100").Value = 2305.5
$ws.Range("J100").Value = 5362.591
$ws.Range("K100").Value = 2305.5
$ws.Range("L100").Value = 5362.591
$ws.Range("M100").Value = -1764.5
$ws.Range("N100").Value = -6444.591

# Row 116
$ws.Range("H116").Value = 3862
$ws.Range("I116").Value = 4281.2
$ws.Range("J116").Value = 3163.3333
$ws.Range("K116").Value = 4281.2
$ws.Range("L116").Value = 3163.3333
$ws.Range("M116").Value = -839.1999999999998
$ws.Range("N116").Value = -10047.3333

# Row 129
$ws.Range("H129").Value = 903.8
$ws.Range("I129").Value = 257.9
$ws.Range("J129").Value = 1047.3334
$ws.Range("K129").Value = 773.6999999999999
$ws.Range("L129").Value = 3142.0002
$ws.Range("M129").Value = 4226.3
$ws.Range("N129").Value = -13142.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 258595.89
$ws.Range("I61").Value = 193786.28
$ws.Range("J61").Value = 388215.06
$ws.Range("K61").Value = 193786.28
$ws.Range("L61").Value = 388215.06
$ws.Range("M61").Value = -193574.28
$ws.Range("N61").Value = -388639.06

# Row 74
$ws.Range("H74").Value = 188038.98
$ws.Range("I74").Value = 244838.7
$ws.Range("J74").Value = 58661.832
$ws.Range("K74").Value = 244838.7
$ws.Range("L74").Value = 58661.832
$ws.Range("M74").Value = -243964.7
$ws.Range("N74").Value = -60409.832

# Row 77
$ws.Range("H77").Value = 188038.98
$ws.Range("I77").Value = 244838.7
$ws.Range("J77").Value = 58661.832
$ws.Range("K77").Value = 1224193.5
$ws.Range("L77").Value = 293309.16
$ws.Range("M77").Value = -1219825.5
$ws.Range("N77").Value = -302045.16

# Row 88
$ws.Range("H88").Value = 2176.6667
$ws.Range("I88").Value = 2081.3794
$ws.Range("J88").Value = 2571.4285
$ws.Range("K88").Value = 2081.3794
$ws.Range("L88").Value = 2571.4285
$ws.Range("M88").Value = -1675.3794
$ws.Range("N88").Value = -3383.4285

# Row 91
$ws.Range("H91").Value = 2176.6667
$ws.Range("I91").Value = 2081.3794
$ws.Range("J91").Value = 2571.4285
$ws.Range("K91").Value = 2081.3794
$ws.Range("L91").Value = 2571.4285
$ws.Range("M91").Value = -677.3793999999998
$ws.Range("N91").Value = -5379.4285

# Row 109
$ws.Range("H109").Value = 33495.75
$ws.Range("J109").Value = 33495.75
$ws.Range("L109").Value = 33495.75
$ws.Range("N109").Value = -36269.75

# Row 132
$ws.Range("H132").Value = 20555.803
$ws.Range("I132").Value = 32012.824
$ws.Range("J132").Value = 2849.5
$ws.Range("K132").Value = 96038.47200000001
$ws.Range("L132").Value = 8548.5
$ws.Range("M132").Value = -93508.47200000001
$ws.Range("N132").Value = -13608.5

# Row 136
$ws.Range("H136").Value = 258595.89
$ws.Range("I136").Value = 193786.28
$ws.Range("J136").Value = 388215.06
$ws.Range("K136").Value = 581358.84
$ws.Range("L136").Value = 1164645.18
$ws.Range("M136").Value = -578808.84
$ws.Range("N136").Value = -1169745.18

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4360
$ws.Range("I86").Value = 5538.5
$ws.Range("J86").Value = 2339.7144
$ws.Range("K86").Value = 5538.5
$ws.Range("L86").Value = 2339.7144
$ws.Range("M86").Value = -4415.5
$ws.Range("N86").Value = -4585.7144

# Row 89
$ws.Range("H89").Value = 4360
$ws.Range("I89").Value = 5538.5
$ws.Range("J89").Value = 2339.7144
$ws.Range("K89").Value = 27692.5
$ws.Range("L89").Value = 11698.572
$ws.Range("M89").Value = -22076.5
$ws.Range("N89").Value = -22930.572

# Row 134
$ws.Range("H134").Value = 7252.8696
$ws.Range("I134").Value = 8361.066000000001
$ws.Range("J134").Value = 5175
$ws.Range("K134").Value = 25083.198
$ws.Range("L134").Value = 15525
$ws.Range("M134").Value = -22548.198
$ws.Range("N134").Value = -20595

$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()

# Row 134
$ws.Range("H134").Value = 1567.8379
$ws.Range("I134").Value = 863.9091
$ws.Range("J134").Value = 2600.2666
$ws.Range("K134").Value = 2591.7273
$ws.Range("L134").Value = 7800.7998
$ws.Range("M134").Value = -56.72730000000001
$ws.Range("N134").Value = -12870.7998

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 70.46154
$ws.Range("I2").Value = 28.2
$ws.Range("J2").Value = 96.875
$ws.Range("K2").Value = 169.2
$ws.Range("L2").Value = 581.25
$ws.Range("M2").Value = -56.19999999999999
$ws.Range("N2").Value = -807.25

# Row 5
$ws.Range("H5").Value = 1067.1613
$ws.Range("I5").Value = 393.9375
$ws.Range("J5").Value = 1785.2667
$ws.Range("K5").Value = 1181.8125
$ws.Range("L5").Value = 5355.800099999999
$ws.Range("M5").Value = -1069.8125
$ws.Range("N5").Value = -5579.800099999999

# Row 38
$ws.Range("H38").Value = 52.77778
$ws.Range("I38").Value = 37.090908
$ws.Range("K38").Value = 111.272724
$ws.Range("M38").Value = 235.727276

# Row 122
$ws.Range("H122").Value = 33333870
$ws.Range("I122").Value = 38461890
$ws.Range("K122").Value = 346157010
$ws.Range("M122").Value = -346154560

# Row 135
$ws.Range("H135").Value = 1067.1613
$ws.Range("I135").Value = 393.9375
$ws.Range("J135").Value = 1785.2667
$ws.Range("K135").Value = 3545.4375
$ws.Range("L135").Value = 16067.4003
$ws.Range("M135").Value = -1010.4375
$ws.Range("N135").Value = -21137.4003

# Row 139
$ws.Range("H139").Value = 2494.5
$ws.Range("I139").Value = 694.26666
$ws.Range("J139").Value = 7895.2
$ws.Range("K139").Value = 2082.79998
$ws.Range("L139").Value = 23685.6
$ws.Range("M139").Value = 3057.20002
$ws.Range("N139").Value = -33965.6

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4638.5264
$ws.Range("I132").Value = 4562.579
$ws.Range("J132").Value = 4714.4736
$ws.Range("K132").Value = 13687.737
$ws.Range("L132").Value = 14143.4208
$ws.Range("M132").Value = -11157.737
$ws.Range("N132").Value = -19203.4208

# Row 136
$ws.Range("H136").Value = 20142
$ws.Range("J136").Value = 20142
$ws.Range("L136").Value = 60426
$ws.Range("N136").Value = -65526

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2900
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 2300
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 2300
$ws.Range("M46").Value = -3312
$ws.Range("N46").Value = -2676

# Row 68
$ws.Range("H68").Value = 2646.5
$ws.Range("I68").Value = 1988.5714
$ws.Range("J68").Value = 3158.2222
$ws.Range("K68").Value = 1988.5714
$ws.Range("L68").Value = 3158.2222
$ws.Range("M68").Value = -1239.5714
$ws.Range("N68").Value = -4656.2222

# Row 71
$ws.Range("H71").Value = 2646.5
$ws.Range("I71").Value = 1988.5714
$ws.Range("J71").Value = 3158.2222
$ws.Range("K71").Value = 9942.857
$ws.Range("L71").Value = 15791.111
$ws.Range("M71").Value = -6198.857
$ws.Range("N71").Value = -23279.111

# Row 86
$ws.Range("H86").Value = 47195
$ws.Range("J86").Value = 47195
$ws.Range("L86").Value = 47195
$ws.Range("N86").Value = -49567

# Row 88
$ws.Range("H88").Value = 29874.75
$ws.Range("J88").Value = 36833
$ws.Range("L88").Value = 36833
$ws.Range("N88").Value = -37689

# Row 89
$ws.Range("H89").Value = 47195
$ws.Range("J89").Value = 47195
$ws.Range("L89").Value = 141585
$ws.Range("N89").Value = -153441

# Row 91
$ws.Range("H91").Value = 29874.75
$ws.Range("J91").Value = 36833
$ws.Range("L91").Value = 36833
$ws.Range("N91").Value = -39797

$ws = $wb.Worksheets.Item("WVR")
# Row 133
$ws.Range("H133").Value = 35235
$ws.Range("J133").Value = 35235
$ws.Range("L133").Value = 35235
$ws.Range("N133").Value = -45355
